# Rerun scripts for 2035 and 2040 EEJ run update
#
# The EEJ (last) run-directory entries for 2035 and 2040 are bumped from
# run "696" to run "697". Setting B19 before B14 keeps the shared-string
# table append order ("2040_06_697" then "2035_06_697") aligned with the
# canonical OOXML produced by Excel itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2040 / EEJ row (row 19): 2040_06_696 -> 2040_06_697
$ws.Range("B19").Value = "2040_06_697"

# 2035 / EEJ row (row 14): 2035_06_696 -> 2035_06_697
$ws.Range("B14").Value = "2035_06_697"

# Leave the active selection on B15, matching the saved view state
$ws.Range("B15").Select()
